$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1200
$ws.Range("I38").Value = 800
$ws.Range("K38").Value = 2400
$ws.Range("M38").Value = -2028
$ws.Range("H39").Value = 333.33334
$ws.Range("I39").Value = 250
$ws.Range("J39").Value = 500
$ws.Range("K39").Value = 750
$ws.Range("L39").Value = 1500
$ws.Range("M39").Value = -454
$ws.Range("N39").Value = -2092
$ws.Range("H76").Value = 9251.5
$ws.Range("I76").Value = 8503
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 8503
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -8188
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 9251.5
$ws.Range("I79").Value = 8503
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 8503
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -7411
$ws.Range("N79").Value = -12184
$ws.Range("H132").Value = 1564.6
$ws.Range("I132").Value = 1476.5
$ws.Range("K132").Value = 4429.5
$ws.Range("M132").Value = -1899.5
$ws.Range("H135").Value = 339.25
$ws.Range("I135").Value = 339.25
$ws.Range("K135").Value = 3053.25
$ws.Range("M135").Value = -518.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("L49").Value = 10000
$ws.Range("N49").Value = -10520
$ws.Range("H61").Value = 2647.1667
$ws.Range("I61").Value = 2261.3333
$ws.Range("J61").Value = 3033
$ws.Range("K61").Value = 2261.3333
$ws.Range("L61").Value = 3033
$ws.Range("M61").Value = -2049.3333
$ws.Range("N61").Value = -3457
$ws.Range("H74").Value = 269015.28
$ws.Range("I74").Value = 309817.7
$ws.Range("K74").Value = 309817.7
$ws.Range("M74").Value = -308943.7
$ws.Range("H77").Value = 269015.28
$ws.Range("I77").Value = 309817.7
$ws.Range("K77").Value = 1549088.5
$ws.Range("M77").Value = -1544720.5
$ws.Range("H136").Value = 2647.1667
$ws.Range("I136").Value = 2261.3333
$ws.Range("J136").Value = 3033
$ws.Range("K136").Value = 6783.999899999999
$ws.Range("L136").Value = 9099
$ws.Range("M136").Value = -4233.999899999999
$ws.Range("N136").Value = -14199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2051.6
$ws.Range("I94").Value = 1613.7778
$ws.Range("J94").Value = 2708.3333
$ws.Range("K94").Value = 1613.7778
$ws.Range("L94").Value = 2708.3333
$ws.Range("M94").Value = -1162.7778
$ws.Range("N94").Value = -3610.3333
$ws.Range("H134").Value = 1991.3334
$ws.Range("I134").Value = 2113.077
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 6339.231000000001
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -3804.231000000001
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 7502.5
$ws.Range("I5").Value = 10007
$ws.Range("J5").Value = 4998
$ws.Range("K5").Value = 10007
$ws.Range("L5").Value = 4998
$ws.Range("M5").Value = -9895
$ws.Range("N5").Value = -5222
$ws.Range("H33").Value = 53498.688
$ws.Range("I33").Value = 41896
$ws.Range("K33").Value = 41896
$ws.Range("M33").Value = -41517
$ws.Range("H38").Value = 15845.833
$ws.Range("J38").Value = 41999.5
$ws.Range("L38").Value = 41999.5
$ws.Range("N38").Value = -42753.5
$ws.Range("H39").Value = 24021.715
$ws.Range("J39").Value = 40500
$ws.Range("L39").Value = 40500
$ws.Range("N39").Value = -41282
$ws.Range("H46").Value = 15845.833
$ws.Range("J46").Value = 41999.5
$ws.Range("L46").Value = 41999.5
$ws.Range("N46").Value = -42421.5
$ws.Range("H49").Value = 24021.715
$ws.Range("J49").Value = 40500
$ws.Range("L49").Value = 40500
$ws.Range("N49").Value = -40864
$ws.Range("H58").Value = 854.5
$ws.Range("I58").Value = 889.3333
$ws.Range("J58").Value = 750
$ws.Range("K58").Value = 889.3333
$ws.Range("L58").Value = 750
$ws.Range("M58").Value = -686.3333
$ws.Range("N58").Value = -1156
$ws.Range("H107").Value = 811.375
$ws.Range("I107").Value = 665.3333
$ws.Range("K107").Value = 665.3333
$ws.Range("M107").Value = 1254.6667
$ws.Range("H132").Value = 2371
$ws.Range("I132").Value = 2371
$ws.Range("K132").Value = 7113
$ws.Range("M132").Value = -4583
$ws.Range("H134").Value = 1987
$ws.Range("I134").Value = 989
$ws.Range("J134").Value = 2985
$ws.Range("K134").Value = 2967
$ws.Range("L134").Value = 8955
$ws.Range("M134").Value = -432
$ws.Range("N134").Value = -14025
$ws.Range("H136").Value = 854.5
$ws.Range("I136").Value = 889.3333
$ws.Range("J136").Value = 750
$ws.Range("K136").Value = 2667.9999
$ws.Range("L136").Value = 2250
$ws.Range("M136").Value = -117.9998999999998
$ws.Range("N136").Value = -7350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1542.2
$ws.Range("J80").Value = 1568.6666
$ws.Range("L80").Value = 1568.6666
$ws.Range("N80").Value = -3564.6666
$ws.Range("H83").Value = 1542.2
$ws.Range("J83").Value = 1568.6666
$ws.Range("L83").Value = 7843.333000000001
$ws.Range("N83").Value = -17827.333
$ws.Range("H107").Value = 1357.3
$ws.Range("I107").Value = 569
$ws.Range("J107").Value = 2145.6
$ws.Range("K107").Value = 569
$ws.Range("L107").Value = 2145.6
$ws.Range("M107").Value = 1351
$ws.Range("N107").Value = -5985.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1879.8
$ws.Range("J22").Value = 2075
$ws.Range("L22").Value = 2075
$ws.Range("N22").Value = -2665
$ws.Range("H27").Value = 1879.8
$ws.Range("J27").Value = 2075
$ws.Range("L27").Value = 2075
$ws.Range("N27").Value = -2289
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 40004490
$ws.Range("I6").Value = 66673668
$ws.Range("J6").Value = 725
$ws.Range("K6").Value = 66673668
$ws.Range("L6").Value = 725
$ws.Range("M6").Value = -66673553
$ws.Range("N6").Value = -955
$ws.Range("H8").Value = 3001.5
$ws.Range("I8").Value = 3001.5
$ws.Range("K8").Value = 3001.5
$ws.Range("M8").Value = -2861.5
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 100000
$ws.Range("L46").Value = 100000
$ws.Range("N46").Value = -100462
$ws.Range("H100").Value = 674.25
$ws.Range("I100").Value = 582.3333
$ws.Range("K100").Value = 1164.6666
$ws.Range("M100").Value = -623.6666
$ws.Range("H113").Value = 257.33334
$ws.Range("I113").Value = 187.5
$ws.Range("J113").Value = 397
$ws.Range("K113").Value = 562.5
$ws.Range("L113").Value = 1191
$ws.Range("M113").Value = 1607.5
$ws.Range("N113").Value = -5531
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 100000
$ws.Range("L134").Value = 300000
$ws.Range("N134").Value = -305070
